$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit ("generated the data for YCbCr color space") fills in six
# GLCM sample rows in column A that were previously left blank/absent in
# the sheet. Excel re-saved the sheet with those rows now present and
# holding a value of 0.
$ws.Range("A41").Value = 0
$ws.Range("A43").Value = 0
$ws.Range("A54").Value = 0
$ws.Range("A55").Value = 0
$ws.Range("A56").Value = 0
$ws.Range("A308").Value = 0

# Match the cursor/scroll position recorded in the saved workbook.
$ws.Application.ActiveWindow.ScrollRow = 385
$ws.Range("A308").Select() | Out-Null
